$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use an existing date-formatted cell as a format template so new date cells
# reuse the workbook's existing date style instead of creating a new numFmt.
$dateTemplate = $ws.Range("A700")

# Row 701
$dateTemplate.Copy()
$ws.Range("A701").PasteSpecial(-4122)
$ws.Range("A701").Value = 45943
$ws.Range("B701").Value = "Kevin Keben"
$ws.Range("C701").Value = "SPECIFIC LINE DEFENSIVE"

# Row 702
$dateTemplate.Copy()
$ws.Range("A702").PasteSpecial(-4122)
$ws.Range("A702").Value = 45943
$ws.Range("B702").Value = "Mattie Pollock"
$ws.Range("C702").Value = "SPECIFIC LINE DEFENSIVE"

# Row 703
$dateTemplate.Copy()
$ws.Range("A703").PasteSpecial(-4122)
$ws.Range("A703").Value = 45943
$ws.Range("B703").Value = "James Abankwah"
$ws.Range("C703").Value = "SPECIFIC LINE DEFENSIVE"

# Row 704
$dateTemplate.Copy()
$ws.Range("A704").PasteSpecial(-4122)
$ws.Range("A704").Value = 45943
$ws.Range("B704").Value = "James Morris"
$ws.Range("C704").Value = "SPECIFIC LINE DEFENSIVE"

# Row 705
$dateTemplate.Copy()
$ws.Range("A705").PasteSpecial(-4122)
$ws.Range("A705").Value = 45943
$ws.Range("B705").Value = "Jeremy Ngakia"
$ws.Range("C705").Value = "SPECIFIC LINE DEFENSIVE"

# Row 706
$dateTemplate.Copy()
$ws.Range("A706").PasteSpecial(-4122)
$ws.Range("A706").Value = 45943
$ws.Range("B706").Value = "Marc Bola"
$ws.Range("C706").Value = "SPECIFIC LINE DEFENSIVE"

# Row 707
$dateTemplate.Copy()
$ws.Range("A707").PasteSpecial(-4122)
$ws.Range("A707").Value = 45943
$ws.Range("B707").Value = "Formose Mendy"
$ws.Range("C707").Value = "SPECIFIC LINE DEFENSIVE"

# Row 708
$dateTemplate.Copy()
$ws.Range("A708").PasteSpecial(-4122)
$ws.Range("A708").Value = 45943
$ws.Range("B708").Value = "Caleb Wiley"
$ws.Range("C708").Value = "SPECIFIC LINE DEFENSIVE"

# Row 709
$dateTemplate.Copy()
$ws.Range("A709").PasteSpecial(-4122)
$ws.Range("A709").Value = 45943
$ws.Range("B709").Value = "Tom Ince"
$ws.Range("C709").Value = "FINISHING PATTERNS"

# Row 710
$dateTemplate.Copy()
$ws.Range("A710").PasteSpecial(-4122)
$ws.Range("A710").Value = 45943
$ws.Range("B710").Value = "Moussa Sissoko"
$ws.Range("C710").Value = "FINISHING PATTERNS"

# Row 711
$dateTemplate.Copy()
$ws.Range("A711").PasteSpecial(-4122)
$ws.Range("A711").Value = 45943
$ws.Range("B711").Value = "Pierre Dwomoh"
$ws.Range("C711").Value = "FINISHING PATTERNS"

# Row 712
$dateTemplate.Copy()
$ws.Range("A712").PasteSpecial(-4122)
$ws.Range("A712").Value = 45943
$ws.Range("B712").Value = "Rocco Vata"
$ws.Range("C712").Value = "FINISHING PATTERNS"

# Row 713
$dateTemplate.Copy()
$ws.Range("A713").PasteSpecial(-4122)
$ws.Range("A713").Value = 45943
$ws.Range("B713").Value = "Luca Kjerrumgaard"
$ws.Range("C713").Value = "FINISHING PATTERNS"

# Row 714
$dateTemplate.Copy()
$ws.Range("A714").PasteSpecial(-4122)
$ws.Range("A714").Value = 45943
$ws.Range("B714").Value = "Kwadwo Baah"
$ws.Range("C714").Value = "FINISHING PATTERNS"

# Row 715
$dateTemplate.Copy()
$ws.Range("A715").PasteSpecial(-4122)
$ws.Range("A715").Value = 45943
$ws.Range("B715").Value = "Amin Nabizada"
$ws.Range("C715").Value = "FINISHING PATTERNS"

# Row 716
$dateTemplate.Copy()
$ws.Range("A716").PasteSpecial(-4122)
$ws.Range("A716").Value = 45943
$ws.Range("B716").Value = "Leo Ramirez Espain"
$ws.Range("C716").Value = "FINISHING PATTERNS"

# Row 717
$dateTemplate.Copy()
$ws.Range("A717").PasteSpecial(-4122)
$ws.Range("A717").Value = 45943
$ws.Range("B717").Value = "Imran Louza"
$ws.Range("C717").Value = "INDIVIDUAL TECHNIQUE - BODY SHAPE & CONTROL/PASS"

# Row 718
$dateTemplate.Copy()
$ws.Range("A718").PasteSpecial(-4122)
$ws.Range("A718").Value = 45944
$ws.Range("B718").Value = "Moussa Sissoko"
$ws.Range("E718").Value = "POSITIONING.PLAY IN THE POCKET"

# Row 719
$dateTemplate.Copy()
$ws.Range("A719").PasteSpecial(-4122)
$ws.Range("A719").Value = 45944
$ws.Range("B719").Value = "Mattie Pollock"
$ws.Range("C719").Value = "INDIVIDUAL SESSION CENTRE BACK (BODY SHAPE IN BUILD UP & RECOVERY PASSES)"

# Row 720
$dateTemplate.Copy()
$ws.Range("A720").PasteSpecial(-4122)
$ws.Range("A720").Value = 45944
$ws.Range("B720").Value = "Formose Mendy"
$ws.Range("C720").Value = "INDIVIDUAL SESSION CENTRE BACK (BODY SHAPE IN BUILD UP & RECOVERY PASSES)"

# Row 721
$dateTemplate.Copy()
$ws.Range("A721").PasteSpecial(-4122)
$ws.Range("A721").Value = 45944
$ws.Range("B721").Value = "Kevin Keben"
$ws.Range("C721").Value = "INDIVIDUAL SESSION CENTRE BACK (BODY SHAPE IN BUILD UP & RECOVERY PASSES)"

# Row 722
$dateTemplate.Copy()
$ws.Range("A722").PasteSpecial(-4122)
$ws.Range("A722").Value = 45944
$ws.Range("B722").Value = "James Abankwah"
$ws.Range("C722").Value = "INDIVIDUAL SESSION CENTRE BACK (BODY SHAPE IN BUILD UP & RECOVERY PASSES)"

$excel.CutCopyMode = 0

# Move the selection to where the user ended up after entering the data.
[void]$ws.Range("D725").Select()
